$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("AA1").Value = "Exp Constant"
$ws.Range("AB1").Value = "Exp Constant [dB]"

$ws.Range("Z1").Copy()
$ws.Range("AA1:AB1").PasteSpecial(-4122)

# New data cells (same values repeated across the three data rows)
$ws.Range("AA2").Value = 385250961.9682089
$ws.Range("AB2").Value = 85.85743731821252

$ws.Range("AA3").Value = 385250961.9682089
$ws.Range("AB3").Value = 85.85743731821252

$ws.Range("AA4").Value = 385250961.9682089
$ws.Range("AB4").Value = 85.85743731821252
